$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 3
$ws.Range("D3").Value = 1.42
$ws.Range("E3").Value = 1.27

# Row 4
$ws.Range("C4").Value = 1.44
$ws.Range("E4").Value = 1.26
$ws.Range("F4").Value = 1.07

# Row 5
$ws.Range("C5").Value = 1.37

# Row 6
$ws.Range("D6").Value = 1.56
$ws.Range("E6").Value = 1.33
$ws.Range("G6").Value = 1.02

# Row 7
$ws.Range("F7").Value = 1.51
